# Booking_Creation_DataSet.xlsx -- refresh the "Created Room ID" column (L)
# with a fresh batch of Room IDs produced by re-running the booking-creation
# test suite. Only the L-column values change; everything else (styles,
# formats, other columns) must stay exactly as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new "Created Room ID" text value
$roomIds = @{
    2 = "173663304"
    9 = "173673069"
    10 = "173674180"
    11 = "173675111"
    12 = "173676578"
    13 = "173677565"
    14 = "173679219"
    15 = "173680239"
    16 = "173681729"
    17 = "173682812"
    18 = "173684339"
    19 = "173685410"
    20 = "173686946"
    21 = "173688180"
    22 = "173689384"
    23 = "173690359"
    24 = "173691683"
    25 = "173692632"
    26 = "173693887"
    27 = "173694863"
    28 = "173696152"
    29 = "173697195"
    30 = "173698631"
    31 = "173699686"
    32 = "173700706"
    33 = "173701697"
    34 = "173702776"
    35 = "173704008"
    36 = "173705476"
    37 = "173706742"
    38 = "173708108"
    39 = "173709307"
    40 = "173710406"
    41 = "173711727"
    42 = "173713154"
    43 = "173714563"
    44 = "173716230"
    45 = "173717623"
    46 = "173719056"
    47 = "173720880"
    48 = "173722930"
    49 = "173724391"
    50 = "173727702"
    51 = "173729324"
    52 = "173731027"
    53 = "173725989"
    54 = "173732735"
    55 = "173734352"
    56 = "173736022"
    57 = "173737528"
    58 = "173739732"
    59 = "173741172"
    60 = "173743335"
    61 = "173770414"
    62 = "173773580"
    63 = "173776590"
    64 = "173779664"
    65 = "173782710"
    66 = "173785874"
    67 = "173789044"
    68 = "173791909"
    69 = "173746659"
    70 = "173760538"
    71 = "173762405"
    72 = "173749294"
    73 = "173753402"
    74 = "173757607"
    75 = "173764376"
    76 = "173766104"
    77 = "173768297"
}

# Scratch cell well outside the used range (column Z, row 1) used only to
# mint a text-typed value; we paste-special VALUES ONLY into the real
# target cell so its existing style/format (s="8" or s="30") is untouched.
$scratch = $ws.Cells.Item(1, 26)

foreach ($row in $roomIds.Keys) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $roomIds[$row]
    $target = $ws.Cells.Item($row, 12)
    $scratch.Copy()
    $target.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

$excel.CutCopyMode = $false
